# Apply the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.323.46'
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("D3").Value = '3.433.51'
$ws.Range("E3").Value = '  +1.73%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.78%  '
$ws.Range("D7").Value = '3.435.59'
$ws.Range("E7").Value = '  +1.78%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.476'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("E11").Value = '  +3.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.387'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.13%  '
$ws.Range("D13").Value = '4.021.97'
$ws.Range("E13").Value = '  +1.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.76%  '
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000173'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.91%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.432.25'
$ws.Range("E17").Value = '  +1.83%  '
$ws.Range("D18").Value = '61.422.39'
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '395.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.564'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").Value = '3.581.71'
$ws.Range("E28").Value = '  +2.10%  '
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("E30").Value = '  +3.48%  '
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("E32").Value = '  -7.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.15'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.00%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.04%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.02'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.88%  '
$ws.Range("B38").Value = 'RenzoRestakedETH'
$ws.Range("C38").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D38").Value = '3.463.95'
$ws.Range("E38").Value = '  +2.13%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.68%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.11'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '167.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0781'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '27.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.802'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.73'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.75%  '
$ws.Range("D49").Value = '2.584.37'
$ws.Range("E49").Value = '  +1.97%  '
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.72%  '
